$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Corridas" to "corridas"
$ws.Name = "corridas"

# Move/update the active cell selection to K24
$ws.Range("K24").Select()
